$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the header labels in row 7 to reflect "Hồ sơ" (dossier) terminology
# instead of "Doanh số" (revenue).
$ws.Range("Q7").Value = "Hồ sơ chi quầy"
$ws.Range("R7").Value = "Hồ sơ chi nhà"

# Edit S7 and T7 text in place first (reusing their existing shared-string
# slots), then swap the two cells' contents so S7 ends up with
# "Hồ sơ chuyển khoản" and T7 ends up with "Tổng".
$ws.Range("S7").Value = "Tổng"
$ws.Range("T7").Value = "Hồ sơ chuyển khoản"

$tmp = $ws.Range("S7").Value2
$ws.Range("S7").Value = $ws.Range("T7").Value2
$ws.Range("T7").Value = $tmp

# Update the active selection to match the edited workbook state
$ws.Range("S7").Select()

$wb.Save()
